# Update "想去人数" (interest count) values in column F across all four sheets
# per the source data refresh (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 2113
$ws.Range("F7").Value = 7949
$ws.Range("F12").Value = 1767
$ws.Range("F13").Value = 1542
$ws.Range("F15").Value = 184
$ws.Range("F16").Value = 4026
$ws.Range("F17").Value = 700
$ws.Range("F19").Value = 1107
$ws.Range("F22").Value = 6220
$ws.Range("F25").Value = 4224
$ws.Range("F26").Value = 705
$ws.Range("F28").Value = 1168
$ws.Range("F34").Value = 49
$ws.Range("F37").Value = 505
$ws.Range("F39").Value = 110
$ws.Range("F41").Value = 153
$ws.Range("F42").Value = 1142

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 372
$ws.Range("F22").Value = 88

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F8").Value = 3099
$ws.Range("F9").Value = 933
$ws.Range("F11").Value = 1265
$ws.Range("F12").Value = 1578

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 3099
$ws.Range("F8").Value = 2113
$ws.Range("F9").Value = 7949
$ws.Range("F11").Value = 933
$ws.Range("F15").Value = 1767
$ws.Range("F16").Value = 1542
$ws.Range("F17").Value = 1265
$ws.Range("F20").Value = 184
$ws.Range("F21").Value = 1578
$ws.Range("F22").Value = 4026
$ws.Range("F23").Value = 372
$ws.Range("F25").Value = 700
$ws.Range("F27").Value = 1107
$ws.Range("F30").Value = 6220
$ws.Range("F32").Value = 705
$ws.Range("F34").Value = 1168
$ws.Range("F39").Value = 88
$ws.Range("F40").Value = 505
$ws.Range("F42").Value = 110
$ws.Range("F44").Value = 1142
